$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; existing rows 32-73 shift down to 33-74.
$ws.Rows.Item(32).Insert()

# Populate the new row 32 with the new record's data.
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C32").Value = 'Arica y Parinacota'
$ws.Range("D32").Value = 44897
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = 100112031
$ws.Range("G32").Value = 'Poroto verde'
$ws.Range("H32").Value = 'Sin especificar'
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 1200
$ws.Range("K32").Value = 650
$ws.Range("L32").Value = 700
$ws.Range("M32").Value = 675
$ws.Range("N32").Value = '$/kilo'
$ws.Range("O32").Value = 'Región de Arica y Parinacota'
$ws.Range("P32").Value = 675
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = 'Hortaliza'
